$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the Test 1 summary label (reuses existing shared string) ---
$ws.Range("C16").Value = "Test 1 Average SUS Score"

# --- 2. New Test 2 header row (row 19) ---
# Order matters for shared-string allocation order, matching original authoring order.
$ws.Range("A19").Value = "Number"
$ws.Range("B19").Value = "User 1 (A)"
$ws.Range("F19").Value = "User 3 (L)"
$ws.Range("H19").Value = "User 4 (JV)"
$ws.Range("D19").Value = "User 2 (JB)"
$ws.Range("C19").Value = "Score Contribution"
$ws.Range("E19").Value = "Score Contribution"
$ws.Range("G19").Value = "Score Contribution"
$ws.Range("I19").Value = "Score Contribution"

# --- 3. Footnote cell next to Test 1 table ---
$ws.Range("G2").Value = "* A SUS score above a 68 would be considered above average and anything below 68 is below average."

# --- 4. Test 2 raw response data (rows 20-29) ---
$numbers = @(1,2,3,4,5,6,7,8,9,10)
$B = @(4,1,5,1,4,2,4,1,5,1)
$D = @(5,1,5,1,5,1,5,1,5,1)
$F = @(3,2,4,1,4,2,5,2,4,2)
$H = @(3,1,5,1,3,3,4,3,5,1)

for ($i = 0; $i -lt 10; $i++) {
    $r = 20 + $i
    $ws.Range("A$r").Value = $numbers[$i]
    $ws.Range("B$r").Value = $B[$i]
    $ws.Range("D$r").Value = $D[$i]
    $ws.Range("F$r").Value = $F[$i]
    $ws.Range("H$r").Value = $H[$i]

    if ($i % 2 -eq 0) {
        $ws.Range("C$r").Formula = "=B$r-1"
        $ws.Range("E$r").Formula = "=D$r-1"
        $ws.Range("G$r").Formula = "=F$r-1"
        $ws.Range("I$r").Formula = "=H$r-1"
    } else {
        $ws.Range("C$r").Formula = "=5-B$r"
        $ws.Range("E$r").Formula = "=5-D$r"
        $ws.Range("G$r").Formula = "=5-F$r"
        $ws.Range("I$r").Formula = "=5-H$r"
    }
}

# --- 5. Totals row (30) ---
$ws.Range("B30").Value = "Total:"
$ws.Range("C30").Formula = "=SUM(C20:C29)"
$ws.Range("D30").Value = "Total:"
$ws.Range("E30").Formula = "=SUM(E20:E29)"
$ws.Range("F30").Value = "Total:"
$ws.Range("G30").Formula = "=SUM(G20:G29)"
$ws.Range("H30").Value = "Total:"
$ws.Range("I30").Formula = "=SUM(I20:I29)"
$ws.Range("B30:M30").Font.Bold = $true

# --- 6. SUS score row (31) ---
$ws.Range("B31").Value = "SUS Score:"
$ws.Range("C31").Formula = "=C30*2.5"
$ws.Range("D31").Value = "SUS Score:"
$ws.Range("E31").Formula = "=E30*2.5"
$ws.Range("F31").Value = "SUS Score:"
$ws.Range("G31").Formula = "=G30*2.5"
$ws.Range("H31").Value = "SUS Score:"
$ws.Range("I31").Formula = "=I30*2.5"
$ws.Range("B31:M31").Font.Bold = $true

# --- 7. Test 2 average SUS score (row 34) ---
$ws.Range("C34").Value = "Test 1 Average SUS Score"
$ws.Range("D34").Formula = "=SUM(C31,E31,G31,I31)/4"
$ws.Range("C34:D34").Font.Bold = $true

# --- 8. Cosmetic: column widths (values tuned so the engine's internal
#        pixel-rounding reproduces the target `width` attribute as closely
#        as representable) ---
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(6).ColumnWidth = 11.750000000000025
$ws.Columns.Item(7).ColumnWidth = 14.416666666666682
$ws.Columns.Item(8).ColumnWidth = 10.083333333333337
$ws.Columns.Item(9).ColumnWidth = 14.416666666666682
$ws.Columns.Item(10).ColumnWidth = 11.583333333333332
$ws.Columns.Item(11).ColumnWidth = 14.416666666666682

# --- 9. Selection / active cell ---
$ws.Range("H9").Select() | Out-Null
